$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "37.469.64"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.074.79"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'235.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6 (XRP)
$ws.Range("D6").Value = "'0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 (Solana)
$ws.Range("D8").Value = "'57.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.395"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.05%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +1.61%  "

# Row 11 (TRON)
$ws.Range("E11").Value = "  +0.98%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "2.378.82"
$ws.Range("E12").Value = "  +0.60%  "

# Row 13 (Chainlink)
$ws.Range("D13").Value = "'14.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "

# Row 14 (Avalanche)
$ws.Range("D14").Value = "'20.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "

# Row 15 (Polygon)
$ws.Range("D15").Value = "'0.781"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "

# Row 16 (Polkadot)
$ws.Range("D16").Value = "'5.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17 (WrappedEther)
$ws.Range("D17").Value = "2.073.20"
$ws.Range("E17").Value = "  +0.54%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "37.391.68"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19 (Uniswap)
$ws.Range("D19").Value = "'6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.96%  "

# Row 20 (Litecoin)
$ws.Range("D20").Value = "'69.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21 (ShibaInu)
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 (BitcoinCash)
$ws.Range("D22").Value = "'227.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23 (Dai)
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 (Toncoin)
$ws.Range("E24").Value = "  +0.88%  "

# Row 25 (PancakeSwap)
$ws.Range("E25").Value = "  -1.99%  "

# Row 26 (Monero)
$ws.Range("D26").Value = "'167.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "

# Row 27 (Cosmos)
$ws.Range("D27").Value = "'8.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "

# Row 28 (ImmutableX)
$ws.Range("E28").Value = "  -1.16%  "

# Row 29 (Kaspa)
$ws.Range("E29").Value = "  +1.25%  "

# Row 30 (EthereumClassic)
$ws.Range("D30").Value = "'19.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31 (Stellar)
$ws.Range("E31").Value = "  -0.52%  "

# Row 32 (Filecoin)
$ws.Range("E32").Value = "  +1.04%  "

# Row 33 (Hedera)
$ws.Range("D33").Value = "'0.0622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.88%  "

# Row 34 (InternetComputer(DFINITY))
$ws.Range("E34").Value = "  +2.46%  "

# Row 35 (LidoDAOToken)
$ws.Range("D35").Value = "'2.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.88%  "

# Row 36 (WEMIXToken)
$ws.Range("E36").Value = "  +1.00%  "

# Row 37 (RenderToken)
$ws.Range("D37").Value = "'3.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

# Row 38 (BinanceUSD)
$ws.Range("E38").Value = "  +0.07%  "

# Row 39 (THORChain)
$ws.Range("D39").Value = "'5.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "

# Row 40: 'HuobiToken' -> 'Cronos'
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "'0.0970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "

# Row 41: 'Aave' -> 'HuobiToken'
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "

# Row 42: 'Cronos' -> 'Aave'
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'98.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.39%  "

# Row 43 (Maker)
$ws.Range("D43").Value = "1.482.43"
$ws.Range("E43").Value = "  +0.11%  "

# Row 44 (TrustWalletToken)
$ws.Range("E44").Value = "  +2.41%  "

# Row 45 (VeChain)
$ws.Range("D45").Value = "'0.0214"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.56%  "

# Row 46 (FTXToken)
$ws.Range("D46").Value = "'4.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.32%  "

# Row 47 (ARBITRUM)
$ws.Range("E47").Value = "  +0.68%  "

# Row 48 (InjectiveProtocol)
$ws.Range("D48").Value = "'15.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.70%  "

# Row 49 (FraxShare)
$ws.Range("D49").Value = "'7.27"
$ws.Range("D49").Style = "Normal"

# Row 50 (MXToken)
$ws.Range("D50").Value = "'2.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "

# Row 51: 'RocketPoolETH' -> 'MultiversX'
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'44.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
